# Update countries & provincias Spain
# Refreshes the COVID-19 country stats table (cols B:H) with the latest
# pull and bumps the "Datos actualizados" timestamp in A1. Because the
# table is kept sorted descending by "Casos totales" (col B), a handful
# of countries with tied/close totals trade row positions as part of
# this refresh - those rows get every column (including the country
# name in col A) rewritten to their new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 21 de Abril de 2020 a las 15:52 -> Datos actualizados a 21 de Abril de 2020 a las 16:22 (table re-sorted by Casos totales)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Abril de 2020 a las 16:22"

# Row 4: Estados Unidos - updated stats
$ws.Cells.Item(4, 2).Value = 794297
$ws.Cells.Item(4, 3).Value = 1538
$ws.Cells.Item(4, 4).Value = 72410
$ws.Cells.Item(4, 5).Value = 679323
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 42564

# Row 8: Alemania - updated stats
$ws.Cells.Item(8, 2).Value = 147593
$ws.Cells.Item(8, 3).Value = 528
$ws.Cells.Item(8, 5).Value = 47531

# Row 17: Paises Bajos - updated stats
$ws.Cells.Item(17, 6).Value = 1087

# Row 60: Moldavia - updated stats
$ws.Cells.Item(60, 2).Value = 2614
$ws.Cells.Item(60, 3).Value = 66
$ws.Cells.Item(60, 5).Value = 2037

# Row 68: Uzbekistan - updated stats
$ws.Cells.Item(68, 4).Value = 344
$ws.Cells.Item(68, 5).Value = 1307

# Row 110: Georgia -> Reunion (table re-sorted by Casos totales)
$ws.Cells.Item(110, 1).Value = "Reunion"
$ws.Cells.Item(110, 2).Value = 410
$ws.Cells.Item(110, 3).Value = 2
$ws.Cells.Item(110, 4).Value = 238
$ws.Cells.Item(110, 5).Value = 172
$ws.Cells.Item(110, 6).Value = 2
$ws.Cells.Item(110, 8).Value = 0

# Row 111: Reunion -> Georgia (table re-sorted by Casos totales)
$ws.Cells.Item(111, 1).Value = "Georgia"
$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 4).Value = 95
$ws.Cells.Item(111, 5).Value = 309
$ws.Cells.Item(111, 6).Value = 6
$ws.Cells.Item(111, 8).Value = 4

# Row 166: Puerto Rico -> Nepal (table re-sorted by Casos totales)
$ws.Cells.Item(166, 1).Value = "Nepal"
$ws.Cells.Item(166, 2).Value = 43
$ws.Cells.Item(166, 3).Value = 12
$ws.Cells.Item(166, 4).Value = 4
$ws.Cells.Item(166, 5).Value = 39
$ws.Cells.Item(166, 8).Value = 0

# Row 167: Eritrea -> Puerto Rico (table re-sorted by Casos totales)
$ws.Cells.Item(167, 1).Value = "Puerto Rico"
$ws.Cells.Item(167, 4).Value = 1
$ws.Cells.Item(167, 8).Value = 2

# Row 168: Siria -> Eritrea (table re-sorted by Casos totales)
$ws.Cells.Item(168, 1).Value = "Eritrea"
$ws.Cells.Item(168, 4).Value = 3
$ws.Cells.Item(168, 5).Value = 36
$ws.Cells.Item(168, 8).Value = 0

# Row 169: Mozambique -> Siria (table re-sorted by Casos totales)
$ws.Cells.Item(169, 1).Value = "Siria"
$ws.Cells.Item(169, 4).Value = 5
$ws.Cells.Item(169, 8).Value = 3

# Row 170: San Martin (Parte Francesa) -> Mozambique (table re-sorted by Casos totales)
$ws.Cells.Item(170, 1).Value = "Mozambique"
$ws.Cells.Item(170, 2).Value = 39
$ws.Cells.Item(170, 4).Value = 8
$ws.Cells.Item(170, 5).Value = 31
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 8).Value = 0

# Row 171: Mongolia -> San Martin (Parte Francesa) (table re-sorted by Casos totales)
$ws.Cells.Item(171, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(171, 2).Value = 37
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 19
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 5
$ws.Cells.Item(171, 8).Value = 2

# Row 172: Republica del Chad -> Mongolia (table re-sorted by Casos totales)
$ws.Cells.Item(172, 1).Value = "Mongolia"
$ws.Cells.Item(172, 2).Value = 34
$ws.Cells.Item(172, 3).Value = 1
$ws.Cells.Item(172, 5).Value = 26

# Row 173: Guam -> Republica del Chad (table re-sorted by Casos totales)
$ws.Cells.Item(173, 1).Value = "Republica del Chad"
$ws.Cells.Item(173, 2).Value = 33
$ws.Cells.Item(173, 4).Value = 8
$ws.Cells.Item(173, 5).Value = 25
$ws.Cells.Item(173, 8).Value = 0

# Row 174: Nepal -> Guam (table re-sorted by Casos totales)
$ws.Cells.Item(174, 1).Value = "Guam"
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 31
$ws.Cells.Item(174, 8).Value = 1

# Row 184: Nueva Caledonia -> Malaui (table re-sorted by Casos totales)
$ws.Cells.Item(184, 1).Value = "Malaui"
$ws.Cells.Item(184, 3).Value = 1
$ws.Cells.Item(184, 4).Value = 3
$ws.Cells.Item(184, 5).Value = 13
$ws.Cells.Item(184, 8).Value = 2

# Row 185: Islas Virgenes de los Estados Unidos -> Nueva Caledonia (table re-sorted by Casos totales)
$ws.Cells.Item(185, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(185, 2).Value = 18
$ws.Cells.Item(185, 4).Value = 17
$ws.Cells.Item(185, 5).Value = 1
$ws.Cells.Item(185, 6).Value = 1

# Row 186: Malaui -> Islas Virgenes de los Estados Unidos (table re-sorted by Casos totales)
$ws.Cells.Item(186, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 17
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 192: Curazao -> Republica de Africa Central (table re-sorted by Casos totales)
$ws.Cells.Item(192, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(192, 3).Value = 2
$ws.Cells.Item(192, 4).Value = 10
$ws.Cells.Item(192, 5).Value = 4
$ws.Cells.Item(192, 8).Value = 0

# Row 193: San Vicente y las Granadinas -> Curazao (table re-sorted by Casos totales)
$ws.Cells.Item(193, 1).Value = "Curazao"
$ws.Cells.Item(193, 2).Value = 14
$ws.Cells.Item(193, 4).Value = 11
$ws.Cells.Item(193, 5).Value = 2
$ws.Cells.Item(193, 8).Value = 1

# Row 194: Republica de Africa Central -> San Vicente y las Granadinas (table re-sorted by Casos totales)
$ws.Cells.Item(194, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(194, 4).Value = 2
$ws.Cells.Item(194, 5).Value = 10

# Row 215: Yemen -> San Pedro y Miquelon (table re-sorted by Casos totales)
$ws.Cells.Item(215, 1).Value = "San Pedro y Miquelon"

# Row 216: San Pedro y Miquelon -> Yemen (table re-sorted by Casos totales)
$ws.Cells.Item(216, 1).Value = "Yemen"
